$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Behandeling")

# Update data values
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 2019

# Update selection on the sheet
$ws.Activate()
$ws.Range("C3").Select()
